$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value as literal TEXT (no numeric auto-conversion),
# then restore the cell style to Normal so no stray number-format style sticks.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = '@'
    $cell.Value = $value
    $cell.Style = 'Normal'
}

# Row 2
Set-TextValue $ws.Range('D2') '66.851.99'
Set-TextValue $ws.Range('E2') '  +0.21%  '

# Row 3
Set-TextValue $ws.Range('D3') '3.879.63'
Set-TextValue $ws.Range('E3') '  +3.97%  '

# Row 4
Set-TextValue $ws.Range('D4') '1.00'
Set-TextValue $ws.Range('E4') '  -0.16%  '

# Row 5
Set-TextValue $ws.Range('D5') '429.41'
Set-TextValue $ws.Range('E5') '  +2.40%  '

# Row 6
Set-TextValue $ws.Range('D6') '131.87'
Set-TextValue $ws.Range('E6') '  +0.23%  '

# Row 7
Set-TextValue $ws.Range('D7') '3.873.65'
Set-TextValue $ws.Range('E7') '  +4.14%  '

# Row 8
Set-TextValue $ws.Range('D8') '0.613'
Set-TextValue $ws.Range('E8') '  -5.73%  '

# Row 10
Set-TextValue $ws.Range('D10') '0.731'

# Row 11
Set-TextValue $ws.Range('D11') '0.167'
Set-TextValue $ws.Range('E11') '  -8.17%  '

# Row 12
Set-TextValue $ws.Range('D12') '0.0000362'
Set-TextValue $ws.Range('E12') '  -10.18%  '

# Row 13
Set-TextValue $ws.Range('D13') '40.91'
Set-TextValue $ws.Range('E13') '  -4.54%  '

# Row 14
Set-TextValue $ws.Range('D14') '4.499.38'
Set-TextValue $ws.Range('E14') '  +4.16%  '

# Row 15
Set-TextValue $ws.Range('D15') '10.07'
Set-TextValue $ws.Range('E15') '  -4.20%  '

# Row 16
Set-TextValue $ws.Range('D16') '15.67'
Set-TextValue $ws.Range('E16') '  +18.06%  '

# Row 17
Set-TextValue $ws.Range('D17') '3.877.42'
Set-TextValue $ws.Range('E17') '  +4.30%  '

# Row 18
Set-TextValue $ws.Range('E18') '  -1.10%  '

# Row 19
Set-TextValue $ws.Range('D19') '19.66'
Set-TextValue $ws.Range('E19') '  -5.35%  '

# Row 20
Set-TextValue $ws.Range('D20') '67.342.62'
Set-TextValue $ws.Range('E20') '  +0.79%  '

# Row 21
Set-TextValue $ws.Range('D21') '1.07'
Set-TextValue $ws.Range('E21') '  -6.07%  '

# Row 22
Set-TextValue $ws.Range('D22') '408.84'
Set-TextValue $ws.Range('E22') '  -8.08%  '

# Row 23
Set-TextValue $ws.Range('D23') '14.46'
Set-TextValue $ws.Range('E23') '  -12.14%  '

# Row 24
Set-TextValue $ws.Range('D24') '85.41'
Set-TextValue $ws.Range('E24') '  -4.90%  '

# Row 25
Set-TextValue $ws.Range('D25') '3.04'
Set-TextValue $ws.Range('E25') '  -4.07%  '

# Row 26
Set-TextValue $ws.Range('D26') '38.10'
Set-TextValue $ws.Range('E26') '  -1.12%  '

# Row 27
Set-TextValue $ws.Range('E27') '  +11.57%  '

# Row 28
Set-TextValue $ws.Range('D28') '3.23'
Set-TextValue $ws.Range('E28') '  -3.15%  '

# Row 29
Set-TextValue $ws.Range('D29') '9.59'
Set-TextValue $ws.Range('E29') '  -6.20%  '

# Row 30
Set-TextValue $ws.Range('D30') '689.91'
Set-TextValue $ws.Range('E30') '  +4.53%  '

# Row 31
Set-TextValue $ws.Range('D31') '0.123'
Set-TextValue $ws.Range('E31') '  -1.50%  '

# Row 32
Set-TextValue $ws.Range('D32') '12.47'
Set-TextValue $ws.Range('E32') '  -2.37%  '

# Row 33
Set-TextValue $ws.Range('E33') '  -0.25%  '

# Row 34
Set-TextValue $ws.Range('D34') '7.16'
Set-TextValue $ws.Range('E34') '  -1.31%  '

# Row 35
Set-TextValue $ws.Range('E35') '  -7.70%  '

# Row 36
Set-TextValue $ws.Range('D36') '38.74'
Set-TextValue $ws.Range('E36') '  -7.69%  '

# Row 37
Set-TextValue $ws.Range('D37') '0.0₃0806'
Set-TextValue $ws.Range('E37') '  +7.44%  '

# Row 38
Set-TextValue $ws.Range('D38') '0.999'
Set-TextValue $ws.Range('E38') '  -0.05%  '

# Row 39
Set-TextValue $ws.Range('D39') '55.32'
Set-TextValue $ws.Range('E39') '  -3.08%  '

# Row 40
$ws.Range('B40').Value = 'ThetaToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextValue $ws.Range('D40') '3.07'
Set-TextValue $ws.Range('E40') '  +0.74%  '

# Row 41
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D41') '0.0458'
Set-TextValue $ws.Range('E41') '  -7.38%  '

# Row 42
Set-TextValue $ws.Range('E42') '  +0.45%  '

# Row 43
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D43') '0.137'
Set-TextValue $ws.Range('E43') '  -9.45%  '

# Row 44
$ws.Range('B44').Value = 'NEARProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D44') '4.52'
Set-TextValue $ws.Range('E44') '  +3.41%  '

# Row 45
Set-TextValue $ws.Range('D45') '147.88'
Set-TextValue $ws.Range('E45') '  +1.02%  '

# Row 46
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D46') '26.81'
Set-TextValue $ws.Range('E46') '  -8.08%  '

# Row 47
Set-TextValue $ws.Range('E47') '  -2.08%  '

# Row 48
$ws.Range('B48').Value = 'LidoDAOToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range('D48') '3.28'
Set-TextValue $ws.Range('E48') '  -5.05%  '

# Row 49
Set-TextValue $ws.Range('E49') '  -4.73%  '

# Row 50
Set-TextValue $ws.Range('D50') '2.78'
Set-TextValue $ws.Range('E50') '  -3.83%  '

# Row 51
Set-TextValue $ws.Range('E51') '  -4.82%  '
